# Auto-generated Excel COM-interop script applying the fr-verb.xlsx diff
# (adds a French/English header row to the existing two sheets and
#  creates four new vocabulary sheets: mouth, leg, hand, sense)
$wb = $excel.ActiveWorkbook
$xlCenter = -4108

# --- Sheet "etre" (#1): add header row (French/English) ---
$wsEtre = $wb.Worksheets.Item(1)
$wsEtre.Range("A1").Value = "French"
$wsEtre.Range("B1").Value = "English"
$wsEtre.Range("A1:B1").HorizontalAlignment = $xlCenter
$wsEtre.Range("A1:B1").VerticalAlignment = $xlCenter
$wsEtre.Activate()
$wsEtre.Range("A1:B1").Select()
$excel.ActiveWindow.Zoom = 145

# --- Sheet "pair" (#2): add header row (French/English x2) ---
$wsPair = $wb.Worksheets.Item(2)
$wsPair.Range("A1").Value = "French"
$wsPair.Range("B1").Value = "English"
$wsPair.Range("C1").Value = "French"
$wsPair.Range("D1").Value = "English"
$wsPair.Range("A1:D1").HorizontalAlignment = $xlCenter
$wsPair.Range("A1:D1").VerticalAlignment = $xlCenter
$wsPair.Activate()
$wsPair.Range("A1:D1").Select()
$excel.ActiveWindow.Zoom = 139

# --- New sheet "mouth" ---
$wsMouth = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMouth.Name = "mouth"
$wsMouth.Columns.Item(1).ColumnWidth = 11.46
$wsMouth.Columns.Item(2).ColumnWidth = 11.69
$wsMouth.Columns.Item(3).ColumnWidth = 9.69
$wsMouth.Columns.Item(4).ColumnWidth = 9.23

$wsMouth.Range("A1").Value = "French"
$wsMouth.Range("B1").Value = "English"
$wsMouth.Range("C1").Value = "French"
$wsMouth.Range("D1").Value = "English"
$wsMouth.Range("E1").Value = "French"
$wsMouth.Range("F1").Value = "English"
$wsMouth.Range("A1:F1").HorizontalAlignment = $xlCenter
$wsMouth.Range("A1:F1").VerticalAlignment = $xlCenter

$wsMouth.Range("A2").Value = "avouer"
$wsMouth.Range("B2").Value = "admit"
$wsMouth.Range("C2").Value = "consentir"
$wsMouth.Range("D2").Value = "consent"
$wsMouth.Range("A2:D2").HorizontalAlignment = $xlCenter
$wsMouth.Range("A2:D2").VerticalAlignment = $xlCenter

$wsMouth.Range("A3").Value = "crier"
$wsMouth.Range("B3").Value = "cry"
$wsMouth.Range("C3").Value = "discuter"
$wsMouth.Range("D3").Value = "discuss"
$wsMouth.Range("A3:D3").HorizontalAlignment = $xlCenter
$wsMouth.Range("A3:D3").VerticalAlignment = $xlCenter

$wsMouth.Range("A4").Value = "nier"
$wsMouth.Range("B4").Value = "deny"
$wsMouth.Range("C4").Value = "exprimer"
$wsMouth.Range("D4").Value = "express"
$wsMouth.Range("A4:D4").HorizontalAlignment = $xlCenter
$wsMouth.Range("A4:D4").VerticalAlignment = $xlCenter

$wsMouth.Range("A5").Value = "rire"
$wsMouth.Range("B5").Value = "laugh"
$wsMouth.Range("C5").Value = "ordonner"
$wsMouth.Range("D5").Value = "order"
$wsMouth.Range("A5:D5").HorizontalAlignment = $xlCenter
$wsMouth.Range("A5:D5").VerticalAlignment = $xlCenter

$wsMouth.Range("A6").Value = "mentir"
$wsMouth.Range("B6").Value = "lie"
$wsMouth.Range("C6").Value = "proposer"
$wsMouth.Range("D6").Value = "propose"
$wsMouth.Range("A6:D6").HorizontalAlignment = $xlCenter
$wsMouth.Range("A6:D6").VerticalAlignment = $xlCenter

$wsMouth.Range("A7").Value = "convaincre"
$wsMouth.Range("B7").Value = "persuade"
$wsMouth.Range("C7").Value = "répondre"
$wsMouth.Range("D7").Value = "response"
$wsMouth.Range("A7:D7").HorizontalAlignment = $xlCenter
$wsMouth.Range("A7:D7").VerticalAlignment = $xlCenter

$wsMouth.Range("A8").Value = "dire"
$wsMouth.Range("B8").Value = "say"
$wsMouth.Range("C8").Value = "suggérer"
$wsMouth.Range("D8").Value = "suggest"
$wsMouth.Range("A8:D8").HorizontalAlignment = $xlCenter
$wsMouth.Range("A8:D8").VerticalAlignment = $xlCenter

$wsMouth.Range("A9").Value = "sourire"
$wsMouth.Range("B9").Value = "smile"
$wsMouth.Range("C9").Value = "confesser"
$wsMouth.Range("D9").Value = "confess"
$wsMouth.Range("A9:D9").HorizontalAlignment = $xlCenter
$wsMouth.Range("A9:D9").VerticalAlignment = $xlCenter

$wsMouth.Range("A10").Value = "parler"
$wsMouth.Range("B10").Value = "speak"
$wsMouth.Range("A10:B10").HorizontalAlignment = $xlCenter
$wsMouth.Range("A10:B10").VerticalAlignment = $xlCenter

$wsMouth.Range("A11").Value = "pleurer"
$wsMouth.Range("B11").Value = "weep"
$wsMouth.Range("A11:B11").HorizontalAlignment = $xlCenter
$wsMouth.Range("A11:B11").VerticalAlignment = $xlCenter

$wsMouth.Range("A12").Value = "manger"
$wsMouth.Range("B12").Value = "eat"
$wsMouth.Range("A12:B12").HorizontalAlignment = $xlCenter
$wsMouth.Range("A12:B12").VerticalAlignment = $xlCenter

$wsMouth.Range("A13").Value = "boire"
$wsMouth.Range("B13").Value = "drink"
$wsMouth.Range("A13:B13").HorizontalAlignment = $xlCenter
$wsMouth.Range("A13:B13").VerticalAlignment = $xlCenter

# placeholder (blank) styled rows left over at the bottom of the sheet
$wsMouth.Range("A14:B14").VerticalAlignment = $xlCenter
$wsMouth.Range("A17:B17").VerticalAlignment = $xlCenter
$wsMouth.Range("A19:B19").VerticalAlignment = $xlCenter
$wsMouth.Range("A20:B20").VerticalAlignment = $xlCenter

$wsMouth.Activate()
$wsMouth.Range("D14").Select()
$excel.ActiveWindow.Zoom = 133

# --- New sheet "leg" ---
$wsLeg = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsLeg.Name = "leg"
$wsLeg.Columns.Item(1).ColumnWidth = 9
$wsLeg.Columns.Item(2).ColumnWidth = 16.38

$wsLeg.Range("A1").Value = "French"
$wsLeg.Range("B1").Value = "English"
$wsLeg.Range("C1").Value = "French"
$wsLeg.Range("D1").Value = "English"
$wsLeg.Range("A1:D1").HorizontalAlignment = $xlCenter
$wsLeg.Range("A1:D1").VerticalAlignment = $xlCenter

$wsLeg.Range("A2").Value = "marcher"
$wsLeg.Range("B2").Value = "walk"
$wsLeg.Range("C2").Value = "danser"
$wsLeg.Range("D2").Value = "dance"
$wsLeg.Range("A2:D2").HorizontalAlignment = $xlCenter
$wsLeg.Range("A2:D2").VerticalAlignment = $xlCenter

$wsLeg.Range("A3").Value = "sauter"
$wsLeg.Range("B3").Value = "jump"
$wsLeg.Range("A3:B3").HorizontalAlignment = $xlCenter
$wsLeg.Range("A3:B3").VerticalAlignment = $xlCenter

$wsLeg.Range("A4").Value = "avancer"
$wsLeg.Range("B4").Value = "marche forward"
$wsLeg.Range("A4:B4").HorizontalAlignment = $xlCenter
$wsLeg.Range("A4:B4").VerticalAlignment = $xlCenter

$wsLeg.Range("A5").Value = "reculer"
$wsLeg.Range("B5").Value = "walk backwards"
$wsLeg.Range("A5:B5").HorizontalAlignment = $xlCenter
$wsLeg.Range("A5:B5").VerticalAlignment = $xlCenter

$wsLeg.Range("A6").Value = "courir"
$wsLeg.Range("B6").Value = "run"
$wsLeg.Range("A6:B6").HorizontalAlignment = $xlCenter
$wsLeg.Range("A6:B6").VerticalAlignment = $xlCenter

# placeholder (blank) styled rows left over at the bottom of the sheet
$wsLeg.Range("A7:B7").VerticalAlignment = $xlCenter
$wsLeg.Range("A8:B8").VerticalAlignment = $xlCenter

$wsLeg.Activate()
$wsLeg.Range("A1:D1").Select()
$excel.ActiveWindow.Zoom = 182

# --- New sheet "hand" ---
$wsHand = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsHand.Name = "hand"

$wsHand.Range("A1").Value = "French"
$wsHand.Range("B1").Value = "English"
$wsHand.Range("C1").Value = "French"
$wsHand.Range("D1").Value = "English"
$wsHand.Range("A1:D1").HorizontalAlignment = $xlCenter
$wsHand.Range("A1:D1").VerticalAlignment = $xlCenter

$wsHand.Activate()
$wsHand.Range("D17").Select()
$excel.ActiveWindow.Zoom = 127

# --- New sheet "sense" ---
$wsSense = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSense.Name = "sense"

$wsSense.Range("A1").Value = "French"
$wsSense.Range("B1").Value = "English"
$wsSense.Range("C1").Value = "French"
$wsSense.Range("D1").Value = "English"
$wsSense.Range("A1:D1").HorizontalAlignment = $xlCenter
$wsSense.Range("A1:D1").VerticalAlignment = $xlCenter

$wsSense.Activate()
$wsSense.Range("A1:D1").Select()
$excel.ActiveWindow.Zoom = 145

# --- Make "mouth" the active tab (activeTab 1 -> 2) ---
$wsMouth.Activate()
$wsMouth.Range("D14").Select()

